# Regenerate the "K" column (column G) of the save_data sheet.
# The previous save used a raw strike-count ("Strike#") for K; this
# regen recalculates K from the underlying pitch log and writes the
# corrected, much smaller per-game values back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value, as produced by the regenerated calc
$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 3
    8  = 2
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 2
    16 = 2
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
